# Hortaliza, Vega Central Mapocho de Santiago - Haba
# Two new weekly records are inserted into the daily price log:
#   - one at row 209 (pushing the former rows 209-231 down by one)
#   - one at row 223 of the resulting sheet (pushing the remainder down by one more)
# ending with the table growing from A1:R231 to A1:R233.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the first new record at row 209 (existing rows 209..231 shift to 210..232)
$ws.Rows("209:209").Insert()

# Insert the second new record at row 223 of the now-shifted sheet
# (existing rows 223..232 shift to 224..233)
$ws.Rows("223:223").Insert()

# --- Fill in the brand-new row 209 ---
$ws.Range("A209").Value = 9
$ws.Range("B209").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C209").Value = "Metropolitana"
$ws.Range("D209").Value = 44748
$ws.Range("E209").Value = 13
$ws.Range("F209").Value = 100112026
$ws.Range("G209").Value = "Haba"
$ws.Range("H209").Value = "Sin especificar"
$ws.Range("I209").Value = "Primera"
$ws.Range("J209").Value = 52
$ws.Range("K209").Value = 20000
$ws.Range("L209").Value = 22000
$ws.Range("M209").Value = 21000
$ws.Range("N209").Value = "`$/saco 25 kilos"
$ws.Range("O209").Value = "Región de Coquimbo"
$ws.Range("P209").Value = 840
$ws.Range("Q209").Value = 25
$ws.Range("R209").Value = "Hortaliza"

# --- Fill in the brand-new row 223 ---
$ws.Range("A223").Value = 9
$ws.Range("B223").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C223").Value = "Metropolitana"
$ws.Range("D223").Value = 44747
$ws.Range("E223").Value = 13
$ws.Range("F223").Value = 100112026
$ws.Range("G223").Value = "Haba"
$ws.Range("H223").Value = "Sin especificar"
$ws.Range("I223").Value = "Primera"
$ws.Range("J223").Value = 52
$ws.Range("K223").Value = 18000
$ws.Range("L223").Value = 18000
$ws.Range("M223").Value = 18000
$ws.Range("N223").Value = "`$/saco 25 kilos"
$ws.Range("O223").Value = "Región de Coquimbo"
$ws.Range("P223").Value = 720
$ws.Range("Q223").Value = 25
$ws.Range("R223").Value = "Hortaliza"
